# Insert a new data row at row 343, shifting the existing rows 343-398
# down to 344-399, then populate the newly inserted row with the new
# price-record data (dated 44522) for "Pepino ensalada".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 343; this pushes the old
# rows 343..398 down to 344..399 and grows the sheet dimension to
# A1:R399 automatically.
$ws.Rows.Item(343).Insert()

# Populate the newly blank row 343 with the new record.
$ws.Range("A343").Value = 8
$ws.Range("B343").Value = "Terminal La Palmera de La Serena"
$ws.Range("C343").Value = "Coquimbo"
$ws.Range("D343").Value = 44522
$ws.Range("E343").Value = 4
$ws.Range("F343").Value = 100112043
$ws.Range("G343").Value = "Pepino ensalada"
$ws.Range("H343").Value = "Sin especificar"
$ws.Range("I343").Value = "Primera"
$ws.Range("J343").Value = 600
$ws.Range("K343").Value = 6500
$ws.Range("L343").Value = 7000
$ws.Range("M343").Value = 6750
$ws.Range("N343").Value = "$/caja 60 unidades"
$ws.Range("O343").Value = "Región de Arica y Parinacota"
$ws.Range("P343").Value = 112
$ws.Range("Q343").Value = 60
$ws.Range("R343").Value = "Hortaliza"
